$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Enterprises (absolute #)" row (12) and the "Enterprises density (per
# 1000 people)" row (13) need to swap places: density now belongs where
# absolute # used to sit, and vice versa. Column A holds the row label,
# column D holds its value.

$absoluteLabel = $ws.Range("A12").Text
$absoluteValue = $ws.Range("D12").Text
$densityLabel  = $ws.Range("A13").Text
$densityValue  = $ws.Range("D13").Text

# Preserve the existing formatting of the value cells so it can be restored
# after the write (see below).
$d12Style = $ws.Range("D12").Style
$d13Style = $ws.Range("D13").Style

# D12/D13 currently contain numeric-looking text ("57227" / "2") that is
# stored as a plain string, not a number. Assigning a numeric-looking string
# straight to .Value would make Excel coerce it into a real number, so the
# cells are temporarily switched to Text format for the assignment and then
# restored to their original style.
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"

$ws.Range("A12").Value = $densityLabel
$ws.Range("D12").Value = $densityValue
$ws.Range("A13").Value = $absoluteLabel
$ws.Range("D13").Value = $absoluteValue

$ws.Range("D12").Style = $d12Style
$ws.Range("D13").Style = $d13Style
